# REPORTE.xlsx update — refresh ESTADO / FECHA DE ULTIMA MODIFICACION values
# for a handful of rows in Hoja1 (sheet1), matching the latest upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 99: I11D71b2d4b -> now PROCESADA, last modified 07/04/2025
$ws.Range("B99").Value = "PROCESADA"
$ws.Range("C99").Value = 45842

# Rows 109-111: I11D71c1c4c, I11D71c1c4d, I11D71c1d3c -> now PENDIENTE (no date)
$ws.Range("B109").Value = "PENDIENTE"
$ws.Range("C109").Value = ""

$ws.Range("B110").Value = "PENDIENTE"
$ws.Range("C110").Value = ""

$ws.Range("B111").Value = "PENDIENTE"
$ws.Range("C111").Value = ""

# Rows 194, 195, 196, 202, 203: were PENDIENTE -> now VOLADA, last modified 07/04/2025
$ws.Range("B194").Value = "VOLADA"
$ws.Range("C194").Value = 45842

$ws.Range("B195").Value = "VOLADA"
$ws.Range("C195").Value = 45842

$ws.Range("B196").Value = "VOLADA"
$ws.Range("C196").Value = 45842

$ws.Range("B202").Value = "VOLADA"
$ws.Range("C202").Value = 45842

$ws.Range("B203").Value = "VOLADA"
$ws.Range("C203").Value = 45842
